$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results for the 380 kV case (Case_3_255, pl_mw)
# Each row (2-25) gets new values in columns B,C,D,F,G,I,J,K,L,O

# Row 2
$ws.Range("B2").Value = 0.6494732456125973
$ws.Range("C2").Value = 0.09309770609986145
$ws.Range("D2").Value = 0.1498256648428224
$ws.Range("F2").Value = 1.821726700776487
$ws.Range("G2").Value = 0.002502836565403285
$ws.Range("I2").Value = 1.1725974267219
$ws.Range("J2").Value = 0.2020730296455824
$ws.Range("K2").Value = 0.3848363136969795
$ws.Range("L2").Value = 0.3572670408846221
$ws.Range("O2").Value = 4.681260454064613

# Row 3
$ws.Range("B3").Value = 0.6089699927461822
$ws.Range("C3").Value = 0.09140758726643838
$ws.Range("D3").Value = 0.1474393341655613
$ws.Range("F3").Value = 1.831029838982673
$ws.Range("G3").Value = 0.002505262397327975
$ws.Range("I3").Value = 1.182315750378546
$ws.Range("J3").Value = 0.2031976675818257
$ws.Range("K3").Value = 0.3477678295008104
$ws.Range("L3").Value = 0.3518618238052369
$ws.Range("O3").Value = 4.713890662532549

# Row 4
$ws.Range("B4").Value = 0.5842533235204428
$ws.Range("C4").Value = 0.09036100408651038
$ws.Range("D4").Value = 0.146024966206987
$ws.Range("F4").Value = 1.83755886458583
$ws.Range("G4").Value = 0.002506832269440347
$ws.Range("I4").Value = 1.188784407837343
$ws.Range("J4").Value = 0.2039568026028391
$ws.Range("K4").Value = 0.3250351446052093
$ws.Range("L4").Value = 0.3486816261462593
$ws.Range("O4").Value = 4.736068753387755

# Row 5
$ws.Range("B5").Value = 0.5742202045479985
$ws.Range("C5").Value = 0.089932312070097
$ws.Range("D5").Value = 0.145461465623363
$ws.Range("F5").Value = 1.840425095327767
$ws.Range("G5").Value = 0.002507492281458638
$ws.Range("I5").Value = 1.191546632966229
$ws.Range("J5").Value = 0.2042834337443935
$ws.Range("K5").Value = 0.3157789270555895
$ws.Range("L5").Value = 0.3474206870474035
$ws.Range("O5").Value = 4.745645682409489

# Row 6
$ws.Range("B6").Value = 0.5725565987764867
$ws.Range("C6").Value = 0.08986099583614759
$ws.Range("D6").Value = 0.1453686761868482
$ws.Range("F6").Value = 1.840913454168295
$ws.Range("G6").Value = 0.002507603102479268
$ws.Range("I6").Value = 1.192012922297494
$ws.Range("J6").Value = 0.2043387148986184
$ws.Range("K6").Value = 0.3142424125586842
$ws.Range("L6").Value = 0.3472134286867004
$ws.Range("O6").Value = 4.747268497952618

# Row 7
$ws.Range("B7").Value = 0.5841178539211285
$ws.Range("C7").Value = 0.0903552314758258
$ws.Range("D7").Value = 0.146017314441778
$ws.Range("F7").Value = 1.837596686854226
$ws.Range("G7").Value = 0.002506841088564012
$ws.Range("I7").Value = 1.188821149119605
$ws.Range("J7").Value = 0.2039611376729162
$ws.Range("K7").Value = 0.3249102806844917
$ws.Range("L7").Value = 0.3486644786863451
$ws.Range("O7").Value = 4.736195727700434

# Row 8
$ws.Range("B8").Value = 0.6354765277864658
$ws.Range("C8").Value = 0.0925168017684328
$ws.Range("D8").Value = 0.1489923510353179
$ws.Range("F8").Value = 1.824764997877452
$ws.Range("G8").Value = 0.002503656340655546
$ws.Range("I8").Value = 1.175844247511716
$ws.Range("J8").Value = 0.2024465822911772
$ws.Range("K8").Value = 0.3720497480207428
$ws.Range("L8").Value = 0.3553746424176296
$ws.Range("O8").Value = 4.692066746156485

# Row 9
$ws.Range("B9").Value = 0.7373717816691681
$ws.Range("C9").Value = 0.09668465262055292
$ws.Range("D9").Value = 0.1552267622091392
$ws.Range("F9").Value = 1.806075813084547
$ws.Range("G9").Value = 0.002498046317756825
$ws.Range("I9").Value = 1.154373169539166
$ws.Range("J9").Value = 0.2000197254608373
$ws.Range("K9").Value = 0.4646859092740101
$ws.Range("L9").Value = 0.3696272214908021
$ws.Range("O9").Value = 4.622521764530944

# Row 10
$ws.Range("B10").Value = 0.8129232214489832
$ws.Range("C10").Value = 0.09970268798851123
$ws.Range("D10").Value = 0.1600477078003451
$ws.Range("F10").Value = 1.796282229133098
$ws.Range("G10").Value = 0.002494308131119838
$ws.Range("I10").Value = 1.141018424225386
$ws.Range("J10").Value = 0.1985663399969049
$ws.Range("K10").Value = 0.5328416289239328
$ws.Range("L10").Value = 0.3807589045761972
$ws.Range("O10").Value = 4.581771284340221

# Row 11
$ws.Range("B11").Value = 0.8474370647994363
$ws.Range("C11").Value = 0.1010659446500384
$ws.Range("D11").Value = 0.1622923921473927
$ws.Range("F11").Value = 1.792679978136761
$ws.Range("G11").Value = 0.002492690004074857
$ws.Range("I11").Value = 1.1354676319513
$ws.Range("J11").Value = 0.1979764167054796
$ws.Range("K11").Value = 0.5638637253270247
$ws.Range("L11").Value = 0.3859650665069978
$ws.Range("O11").Value = 4.565476097132802

# Row 12
$ws.Range("B12").Value = 0.8605267454438206
$ws.Range("C12").Value = 0.1015807659267836
$ws.Range("D12").Value = 0.1631497449522357
$ws.Range("F12").Value = 1.791438377155885
$ws.Range("G12").Value = 0.002492089049019129
$ws.Range("I12").Value = 1.133441029035204
$ws.Range("J12").Value = 0.1977632450465983
$ws.Range("K12").Value = 0.5756130019905186
$ws.Range("L12").Value = 0.3879568273357137
$ws.Range("O12").Value = 4.559627747967653

# Row 13
$ws.Range("B13").Value = 0.8577067705442118
$ws.Range("C13").Value = 0.1014699532563128
$ws.Range("D13").Value = 0.1629647735649655
$ws.Range("F13").Value = 1.791700332653392
$ws.Range("G13").Value = 0.002492217951615877
$ws.Range("I13").Value = 1.133874143123542
$ws.Range("J13").Value = 0.1978087012311676
$ws.Range("K13").Value = 0.5730825119155156
$ws.Range("L13").Value = 0.3875269657343523
$ws.Range("O13").Value = 4.560872963637905

# Row 14
$ws.Range("B14").Value = 0.8485135638742065
$ws.Range("C14").Value = 0.1011083277947264
$ws.Range("D14").Value = 0.1623627804534067
$ws.Range("F14").Value = 1.792575376803896
$ws.Range("G14").Value = 0.002492640327053686
$ws.Range("I14").Value = 1.135299391972843
$ws.Range("J14").Value = 0.1979586742636137
$ws.Range("K14").Value = 0.5648303114774649
$ws.Range("L14").Value = 0.3861285239534453
$ws.Range("O14").Value = 4.564988491846577

# Row 15
$ws.Range("B15").Value = 0.842885044182168
$ws.Range("C15").Value = 0.1008866368197729
$ws.Range("D15").Value = 0.1619949956809137
$ws.Range("F15").Value = 1.79312731432104
$ws.Range("G15").Value = 0.002492900578696853
$ws.Range("I15").Value = 1.136182211202851
$ws.Range("J15").Value = 0.198051867269065
$ws.Range("K15").Value = 0.5597758258012107
$ws.Range("L15").Value = 0.3852745768512449
$ws.Range("O15").Value = 4.567551337585172

# Row 16
$ws.Range("B16").Value = 0.8106705096805911
$ws.Range("C16").Value = 0.09961339938370628
$ws.Range("D16").Value = 0.1599020439366825
$ws.Range("F16").Value = 1.79653479418549
$ws.Range("G16").Value = 0.002494415533013747
$ws.Range("I16").Value = 1.141391729446052
$ws.Range("J16").Value = 0.1986063242001066
$ws.Range("K16").Value = 0.5308145615354363
$ws.Range("L16").Value = 0.3804215192441802
$ws.Range("O16").Value = 4.582881321538082

# Row 17
$ws.Range("B17").Value = 0.7909445121383385
$ws.Range("C17").Value = 0.09882981679305658
$ws.Range("D17").Value = 0.1586312463934547
$ws.Range("F17").Value = 1.798843512225936
$ws.Range("G17").Value = 0.00249536597308174
$ws.Range("I17").Value = 1.144721870538149
$ws.Range("J17").Value = 0.1989646925326376
$ws.Range("K17").Value = 0.5130518275572911
$ws.Range("L17").Value = 0.3774806570845328
$ws.Range("O17").Value = 4.592859983448761

# Row 18
$ws.Range("B18").Value = 0.7796123484140765
$ws.Range("C18").Value = 0.09837821278733827
$ws.Range("D18").Value = 0.1579051801316922
$ws.Range("F18").Value = 1.800251713850955
$ws.Range("G18").Value = 0.002495920399228932
$ws.Range("I18").Value = 1.146686636081984
$ws.Range("J18").Value = 0.1991775223145567
$ws.Range("K18").Value = 0.5028368801218619
$ws.Range("L18").Value = 0.3758025549556692
$ws.Range("O18").Value = 4.59881052239254

# Row 19
$ws.Range("B19").Value = 0.7757778522403385
$ws.Range("C19").Value = 0.0982251522385198
$ws.Range("D19").Value = 0.157660184280715
$ws.Range("F19").Value = 1.800742300926437
$ws.Range("G19").Value = 0.002496109452660916
$ws.Range("I19").Value = 1.147360350149366
$ws.Range("J19").Value = 0.1992507352558093
$ws.Range("K19").Value = 0.4993785879935615
$ws.Range("L19").Value = 0.3752366855198801
$ws.Range("O19").Value = 4.600861530459269

# Row 20
$ws.Range("B20").Value = 0.7930429655115745
$ws.Range("C20").Value = 0.0989133246567846
$ws.Range("D20").Value = 0.1587660222684519
$ws.Range("F20").Value = 1.798589436838938
$ws.Range("G20").Value = 0.002495263994398514
$ws.Range("I20").Value = 1.144362263407007
$ws.Range("J20").Value = 0.1989258497553159
$ws.Range("K20").Value = 0.5149425290380805
$ws.Range("L20").Value = 0.3777923307974334
$ws.Range("O20").Value = 4.591775893460294

# Row 21
$ws.Range("B21").Value = 0.8512132943194786
$ws.Range("C21").Value = 0.1012145845444223
$ws.Range("D21").Value = 0.1625394018839756
$ws.Range("F21").Value = 1.79231503191032
$ws.Range("G21").Value = 0.002492515945594289
$ws.Range("I21").Value = 1.134878716853727
$ws.Range("J21").Value = 0.1979143463789512
$ws.Range("K21").Value = 0.5672541370007593
$ws.Range("L21").Value = 0.3865387305557988
$ws.Range("O21").Value = 4.563770916105028

# Row 22
$ws.Range("B22").Value = 0.8893473755527452
$ws.Range("C22").Value = 0.1027103350444349
$ws.Range("D22").Value = 0.1650482660985375
$ws.Range("F22").Value = 1.788928259218409
$ws.Range("G22").Value = 0.002490788659286863
$ws.Range("I22").Value = 1.129119918346973
$ws.Range("J22").Value = 0.1973128267504478
$ws.Range("K22").Value = 0.6014534591332108
$ws.Range("L22").Value = 0.3923732585700321
$ws.Range("O22").Value = 4.547346486041278

# Row 23
$ws.Range("B23").Value = 0.8689841209834412
$ws.Range("C23").Value = 0.1019127884621724
$ws.Range("D23").Value = 0.1637053542332296
$ws.Range("F23").Value = 1.790670570892203
$ws.Range("G23").Value = 0.002491704273769555
$ws.Range("I23").Value = 1.132153319169596
$ws.Range("J23").Value = 0.1976284273252986
$ws.Range("K23").Value = 0.5831998811445942
$ws.Range("L23").Value = 0.3892484950274309
$ws.Range("O23").Value = 4.555940692601808

# Row 24
$ws.Range("B24").Value = 0.7920942279685619
$ws.Range("C24").Value = 0.09887557421250648
$ws.Range("D24").Value = 0.1587050759721222
$ws.Range("F24").Value = 1.798704052346281
$ws.Range("G24").Value = 0.002495310074148326
$ws.Range("I24").Value = 1.144524685345658
$ws.Range("J24").Value = 0.1989433893951791
$ws.Range("K24").Value = 0.5140877519294236
$ws.Range("L24").Value = 0.3776513837324131
$ws.Range("O24").Value = 4.592265344949425

# Row 25
$ws.Range("B25").Value = 0.7096831907826129
$ws.Range("C25").Value = 0.09556482223538865
$ws.Range("D25").Value = 0.1534976940453419
$ws.Range("F25").Value = 1.810439626686708
$ws.Range("G25").Value = 0.002499496360575199
$ws.Range("I25").Value = 1.159756377057903
$ws.Range("J25").Value = 0.2006182597905166
$ws.Range("K25").Value = 0.4396068417128731
$ws.Range("L25").Value = 0.3656550870597073
$ws.Range("O25").Value = 4.639517929082672
